$d = $word.ActiveDocument

$ids = @("p054v_1", "p054v_2", "p054v_3", "p054v_4", "p054v_5")

foreach ($id in $ids) {
    $find = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $find, 2)
}
